# Natmi following Dr Hou advice:
# A new "ECs" target cluster is introduced and the FAPs/sCs x FAPs/sCs/ECs
# sending/target combinations (and their NATMI statistics) are refreshed
# to reflect the updated 3-sample run (rows 2-7 of Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgf16"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.678104
$ws.Range("H2").Value = 5.034312
$ws.Range("I2").Value = 0.551436927751233
$ws.Range("J2").Value = 0.551436927751233
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 4.220261333333333
$ws.Range("N2").Value = 12.660784
$ws.Range("O2").Value = 0.6739259863235564
$ws.Range("P2").Value = 0.6739259863235564
$ws.Range("Q2").Value = 7.082037424512
$ws.Range("R2").Value = 63.738336820608
$ws.Range("S2").Value = 0.3716276754299814
$ws.Range("T2").Value = 0.3716276754299814

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgf16"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.678104
$ws.Range("H3").Value = 5.034312
$ws.Range("I3").Value = 0.551436927751233
$ws.Range("J3").Value = 0.551436927751233
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.296447666666667
$ws.Range("N3").Value = 3.889343
$ws.Range("O3").Value = 0.2070274097896007
$ws.Range("P3").Value = 0.2070274097896007
$ws.Range("Q3").Value = 2.175574015224
$ws.Range("R3").Value = 19.580166137016
$ws.Range("S3").Value = 0.1141625588146729
$ws.Range("T3").Value = 0.1141625588146729

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf16"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.678104
$ws.Range("H4").Value = 5.034312
$ws.Range("I4").Value = 0.551436927751233
$ws.Range("J4").Value = 0.551436927751233
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.745494
$ws.Range("N4").Value = 2.236482
$ws.Range("O4").Value = 0.119046603886843
$ws.Range("P4").Value = 0.119046603886843
$ws.Range("Q4").Value = 1.251016463376
$ws.Range("R4").Value = 11.259148170384
$ws.Range("S4").Value = 0.06564669350657872
$ws.Range("T4").Value = 0.06564669350657872

$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Fgf16"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.365043666666667
$ws.Range("H5").Value = 4.095131
$ws.Range("I5").Value = 0.448563072248767
$ws.Range("J5").Value = 0.448563072248767
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 4.220261333333333
$ws.Range("N5").Value = 12.660784
$ws.Range("O5").Value = 0.6739259863235564
$ws.Range("P5").Value = 0.6739259863235564
$ws.Range("Q5").Value = 5.760841004744889
$ws.Range("R5").Value = 51.847569042704
$ws.Range("S5").Value = 0.302298310893575
$ws.Range("T5").Value = 0.302298310893575

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Fgf16"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.365043666666667
$ws.Range("H6").Value = 4.095131
$ws.Range("I6").Value = 0.448563072248767
$ws.Range("J6").Value = 0.448563072248767
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.296447666666667
$ws.Range("N6").Value = 3.889343
$ws.Range("O6").Value = 0.2070274097896007
$ws.Range("P6").Value = 0.2070274097896007
$ws.Range("Q6").Value = 1.769707676548111
$ws.Range("R6").Value = 15.927369088933
$ws.Range("S6").Value = 0.09286485097492773
$ws.Range("T6").Value = 0.09286485097492775

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Fgf16"
$ws.Range("C7").Value = "Fgfr3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.365043666666667
$ws.Range("H7").Value = 4.095131
$ws.Range("I7").Value = 0.448563072248767
$ws.Range("J7").Value = 0.448563072248767
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.745494
$ws.Range("N7").Value = 2.236482
$ws.Range("O7").Value = 0.119046603886843
$ws.Range("P7").Value = 0.119046603886843
$ws.Range("Q7").Value = 1.017631863238
$ws.Range("R7").Value = 9.158686769142001
$ws.Range("S7").Value = 0.05339991038026431
$ws.Range("T7").Value = 0.05339991038026431
